# Natmi following Dr Hou advice
# Update the LR-pairs matrix for Myoc-Fzd3 to include the "ECs" cluster,
# expanding the Sending cluster x Target cluster combinations from a
# 2x2 grid (FAPs/sCs) to a full 3x3 grid (ECs/FAPs/sCs), with refreshed
# statistics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Myoc"
$ws.Range("C2").Value = "Fzd3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2144083333333333
$ws.Range("H2").Value = 0.643225
$ws.Range("I2").Value = 0.008611346839948651
$ws.Range("J2").Value = 0.008611346839948651
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.173174
$ws.Range("N2").Value = 0.519522
$ws.Range("O2").Value = 0.0473074116693291
$ws.Range("P2").Value = 0.0473074116693291
$ws.Range("Q2").Value = 0.03712994871666667
$ws.Range("R2").Value = 0.33416953845
$ws.Range("S2").Value = 0.0004073805299848271
$ws.Range("T2").Value = 0.0004073805299848271

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Myoc"
$ws.Range("C3").Value = "Fzd3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2144083333333333
$ws.Range("H3").Value = 0.643225
$ws.Range("I3").Value = 0.008611346839948651
$ws.Range("J3").Value = 0.008611346839948651
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.880936
$ws.Range("N3").Value = 2.642808
$ws.Range("O3").Value = 0.2406527654632456
$ws.Range("P3").Value = 0.2406527654632456
$ws.Range("Q3").Value = 0.1888800195333334
$ws.Range("R3").Value = 1.6999201758
$ws.Range("S3").Value = 0.002072344431396824
$ws.Range("T3").Value = 0.002072344431396824

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Myoc"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2144083333333333
$ws.Range("H4").Value = 0.643225
$ws.Range("I4").Value = 0.008611346839948651
$ws.Range("J4").Value = 0.008611346839948651
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.606500333333333
$ws.Range("N4").Value = 7.819501
$ws.Range("O4").Value = 0.7120398228674253
$ws.Range("P4").Value = 0.7120398228674253
$ws.Range("Q4").Value = 0.5588553923027778
$ws.Range("R4").Value = 5.029698530725001
$ws.Range("S4").Value = 0.006131621878567
$ws.Range("T4").Value = 0.006131621878567

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Myoc"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.971258
$ws.Range("H5").Value = 71.91377399999999
$ws.Range("I5").Value = 0.9627648963950115
$ws.Range("J5").Value = 0.9627648963950115
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.173174
$ws.Range("N5").Value = 0.519522
$ws.Range("O5").Value = 0.0473074116693291
$ws.Range("P5").Value = 0.0473074116693291
$ws.Range("Q5").Value = 4.151198632892
$ws.Range("R5").Value = 37.360787696028
$ws.Range("S5").Value = 0.04554591529453779
$ws.Range("T5").Value = 0.04554591529453779

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Myoc"
$ws.Range("C6").Value = "Fzd3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 23.971258
$ws.Range("H6").Value = 71.91377399999999
$ws.Range("I6").Value = 0.9627648963950115
$ws.Range("J6").Value = 0.9627648963950115
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.880936
$ws.Range("N6").Value = 2.642808
$ws.Range("O6").Value = 0.2406527654632456
$ws.Range("P6").Value = 0.2406527654632456
$ws.Range("Q6").Value = 21.117144137488
$ws.Range("R6").Value = 190.054297237392
$ws.Range("S6").Value = 0.2316920348083946
$ws.Range("T6").Value = 0.2316920348083947

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Myoc"
$ws.Range("C7").Value = "Fzd3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 23.971258
$ws.Range("H7").Value = 71.91377399999999
$ws.Range("I7").Value = 0.9627648963950115
$ws.Range("J7").Value = 0.9627648963950115
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.606500333333333
$ws.Range("N7").Value = 7.819501
$ws.Range("O7").Value = 0.7120398228674253
$ws.Range("P7").Value = 0.7120398228674253
$ws.Range("Q7").Value = 62.48109196741932
$ws.Range("R7").Value = 562.3298277067739
$ws.Range("S7").Value = 0.685526946292079
$ws.Range("T7").Value = 0.685526946292079

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Myoc"
$ws.Range("C8").Value = "Fzd3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7126843333333334
$ws.Range("H8").Value = 2.138053
$ws.Range("I8").Value = 0.02862375676503981
$ws.Range("J8").Value = 0.02862375676503981
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.173174
$ws.Range("N8").Value = 0.519522
$ws.Range("O8").Value = 0.0473074116693291
$ws.Range("P8").Value = 0.0473074116693291
$ws.Range("Q8").Value = 0.1234183967406667
$ws.Range("R8").Value = 1.110765570666
$ws.Range("S8").Value = 0.001354115844806482
$ws.Range("T8").Value = 0.001354115844806482

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Myoc"
$ws.Range("C9").Value = "Fzd3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7126843333333334
$ws.Range("H9").Value = 2.138053
$ws.Range("I9").Value = 0.02862375676503981
$ws.Range("J9").Value = 0.02862375676503981
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.880936
$ws.Range("N9").Value = 2.642808
$ws.Range("O9").Value = 0.2406527654632456
$ws.Range("P9").Value = 0.2406527654632456
$ws.Range("Q9").Value = 0.6278292858693334
$ws.Range("R9").Value = 5.650463572824001
$ws.Range("S9").Value = 0.006888386223454115
$ws.Range("T9").Value = 0.006888386223454116

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Myoc"
$ws.Range("C10").Value = "Fzd3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7126843333333334
$ws.Range("H10").Value = 2.138053
$ws.Range("I10").Value = 0.02862375676503981
$ws.Range("J10").Value = 0.02862375676503981
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.606500333333333
$ws.Range("N10").Value = 7.819501
$ws.Range("O10").Value = 0.7120398228674253
$ws.Range("P10").Value = 0.7120398228674253
$ws.Range("Q10").Value = 1.857611952394778
$ws.Range("R10").Value = 16.718507571553
$ws.Range("S10").Value = 0.02038125469677922
$ws.Range("T10").Value = 0.02038125469677922

